# Update header metadata (timestamp + row counts) on all three sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 17:36:10"
$ws.Range("A3").Value = "Total filas: 285"

$ws.Range("A23").Value = "05:57:13"
$ws.Range("B23").Value = "07:21"
$ws.Range("C23").Value = "23_HERNANDEZ"
$ws.Range("D23").Value = 84
$ws.Range("E23").Value = "LP1912"
$ws.Range("A24").Value = "06:17:28"
$ws.Range("B24").Value = "07:21"
$ws.Range("C24").Value = "16_SANTA ANA"
$ws.Range("D24").Value = 64
$ws.Range("E24").Value = "LP1912"
$ws.Range("A40").Value = "06:35:22"
$ws.Range("B40").Value = "08:29"
$ws.Range("C40").Value = "11_ETCHEVERRY"
$ws.Range("D40").Value = 114
$ws.Range("E40").Value = "LP1912"
$ws.Range("A41").Value = "06:35:22"
$ws.Range("B41").Value = "08:29"
$ws.Range("C41").Value = "15_ABASTO"
$ws.Range("D41").Value = 114
$ws.Range("E41").Value = "LP1912"
$ws.Range("A58").Value = "07:38:39"
$ws.Range("B58").Value = "09:18"
$ws.Range("C58").Value = "15X38_ABASTO"
$ws.Range("D58").Value = 100
$ws.Range("E58").Value = "LP1912"
$ws.Range("A59").Value = "08:10:18"
$ws.Range("B59").Value = "09:18"
$ws.Range("C59").Value = "14_ABASTO"
$ws.Range("D59").Value = 68
$ws.Range("E59").Value = "LP1912"
$ws.Range("A111").Value = "10:50:41"
$ws.Range("B111").Value = "11:54"
$ws.Range("C111").Value = "23_HERNANDEZ"
$ws.Range("D111").Value = 64
$ws.Range("E111").Value = "LP1912"
$ws.Range("A113").Value = "11:54:18"
$ws.Range("B113").Value = "11:54"
$ws.Range("C113").Value = "15X38_ABASTO"
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = "LP1912"
$ws.Range("A120").Value = "10:37:52"
$ws.Range("B120").Value = "12:10"
$ws.Range("C120").Value = "16_P MOR-SANTA ANA"
$ws.Range("D120").Value = 93
$ws.Range("E120").Value = "LP1912"
$ws.Range("A121").Value = "10:37:52"
$ws.Range("B121").Value = "12:10"
$ws.Range("C121").Value = "15_ABASTO"
$ws.Range("D121").Value = 93
$ws.Range("E121").Value = "LP1912"
$ws.Range("A128").Value = "11:47:17"
$ws.Range("B128").Value = "12:32"
$ws.Range("C128").Value = "23_HERNANDEZ"
$ws.Range("D128").Value = 45
$ws.Range("E128").Value = "LP1912"
$ws.Range("A129").Value = "10:37:52"
$ws.Range("B129").Value = "12:32"
$ws.Range("C129").Value = "14_ABASTO"
$ws.Range("D129").Value = 115
$ws.Range("E129").Value = "LP1912"
$ws.Range("A130").Value = "11:34:59"
$ws.Range("B130").Value = "12:33"
$ws.Range("C130").Value = "15_ABASTO"
$ws.Range("D130").Value = 59
$ws.Range("E130").Value = "LP1912"
$ws.Range("A131").Value = "11:47:17"
$ws.Range("B131").Value = "12:33"
$ws.Range("C131").Value = "14_ABASTO"
$ws.Range("D131").Value = 46
$ws.Range("E131").Value = "LP1912"
$ws.Range("A209").Value = "13:56:11"
$ws.Range("B209").Value = "15:53"
$ws.Range("C209").Value = "15X38_ABASTO"
$ws.Range("D209").Value = 117
$ws.Range("E209").Value = "LP1912"
$ws.Range("A211").Value = "13:56:11"
$ws.Range("B211").Value = "15:53"
$ws.Range("C211").Value = "16_P MOR-SANTA ANA"
$ws.Range("D211").Value = 117
$ws.Range("E211").Value = "LP1912"
$ws.Range("A228").Value = "15:46:07"
$ws.Range("B228").Value = "16:30"
$ws.Range("C228").Value = "14_ABASTO"
$ws.Range("D228").Value = 44
$ws.Range("E228").Value = "LP1912"
$ws.Range("A229").Value = "15:17:33"
$ws.Range("B229").Value = "16:30"
$ws.Range("C229").Value = "16_SANTA ANA"
$ws.Range("D229").Value = 73
$ws.Range("E229").Value = "LP1912"
$ws.Range("A279").Value = "17:36:10"
$ws.Range("B279").Value = "18:37"
$ws.Range("C279").Value = "23_HERNANDEZ"
$ws.Range("D279").Value = 61
$ws.Range("E279").Value = "LP1912"
$ws.Range("A280").Value = "17:13:39"
$ws.Range("B280").Value = "18:41"
$ws.Range("C280").Value = "10_OLMOS"
$ws.Range("D280").Value = 88
$ws.Range("E280").Value = "LP1912"
$ws.Range("A281").Value = "16:52:42"
$ws.Range("B281").Value = "18:45"
$ws.Range("C281").Value = "16_SANTA ANA"
$ws.Range("D281").Value = 113
$ws.Range("E281").Value = "LP1912"
$ws.Range("A282").Value = "17:13:39"
$ws.Range("B282").Value = "18:52"
$ws.Range("C282").Value = "17_ROMERO"
$ws.Range("D282").Value = 99
$ws.Range("E282").Value = "LP1912"
$ws.Range("A283").Value = "17:13:39"
$ws.Range("B283").Value = "18:57"
$ws.Range("C283").Value = "16_P MOR-SANTA ANA"
$ws.Range("D283").Value = 104
$ws.Range("E283").Value = "LP1912"
$ws.Range("A284").Value = "17:13:39"
$ws.Range("B284").Value = "18:59"
$ws.Range("C284").Value = "14_ABASTO"
$ws.Range("D284").Value = 106
$ws.Range("E284").Value = "LP1912"
$ws.Range("A285").Value = "17:36:10"
$ws.Range("B285").Value = "19:00"
$ws.Range("C285").Value = "14_ABASTO"
$ws.Range("D285").Value = 84
$ws.Range("E285").Value = "LP1912"
$ws.Range("A286").Value = "17:13:39"
$ws.Range("B286").Value = "19:03"
$ws.Range("C286").Value = "215_EL PELIGRO"
$ws.Range("D286").Value = 110
$ws.Range("E286").Value = "LP1912"
$ws.Range("A287").Value = "17:36:10"
$ws.Range("B287").Value = "19:04"
$ws.Range("C287").Value = "215_EL PELIGRO"
$ws.Range("D287").Value = 88
$ws.Range("E287").Value = "LP1912"
$ws.Range("A288").Value = "17:36:10"
$ws.Range("B288").Value = "19:17"
$ws.Range("C288").Value = "27_EL RETIRO"
$ws.Range("D288").Value = 101
$ws.Range("E288").Value = "LP1912"
$ws.Range("A289").Value = "17:36:10"
$ws.Range("B289").Value = "19:17"
$ws.Range("C289").Value = "14X44_ABASTO"
$ws.Range("D289").Value = 101
$ws.Range("E289").Value = "LP1912"
$ws.Range("A290").Value = "17:36:10"
$ws.Range("B290").Value = "19:28"
$ws.Range("C290").Value = "215C_EL PATO"
$ws.Range("D290").Value = 112
$ws.Range("E290").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 17:36:10"
$ws.Range("A3").Value = "Total filas: 46"

$ws.Range("A50").Value = "17:36:10"
$ws.Range("B50").Value = "19:04"
$ws.Range("C50").Value = "215_EL PELIGRO"
$ws.Range("D50").Value = 88
$ws.Range("E50").Value = "LP1912"
$ws.Range("A51").Value = "17:36:10"
$ws.Range("B51").Value = "19:28"
$ws.Range("C51").Value = "215C_EL PATO"
$ws.Range("D51").Value = 112
$ws.Range("E51").Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 17:36:10"
$ws.Range("A3").Value = "Total filas: 38"

$ws.Range("A43").Value = "17:36:10"
$ws.Range("B43").Value = "19:24"
$ws.Range("C43").Value = "215B_LP-P MOR-1 Y 57"
$ws.Range("D43").Value = 108
$ws.Range("E43").Value = "L6173"
